$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UML diagram header labels (row 4): "Modely ..." -> "Diagramy ..."
$ws.Range("L4").Value = "Diagramy balíčkov (Package)"
$ws.Range("M4").Value = "Diagramy tried (Class)"
$ws.Range("N4").Value = "Sekvenčné diagramy (Sequence)"

# Adjust column M width to fit the new (longer) text
$ws.Columns.Item(13).ColumnWidth = 8.5

# Move the active cell selection as recorded in the saved file
$ws.Range("O13").Select()
